$d = $word.ActiveDocument
$app = $word
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertParagraphAfter()

$p1 = $d.Paragraphs.Last
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last
$p2.Range.Text = "placeholder"
$p2.Style = "List Paragraph"

$gallery = $app.ListGalleries.Item(1)
$tmpl = $gallery.ListTemplates.Item(1)
$p2.Range.ListFormat.ApplyListTemplate($tmpl)

$st = $d.Styles.Item("List Paragraph")
$st.Priority = 34
$st.ParagraphFormat.LeftIndent = 36
$st.NoSpaceBetweenParagraphsOfSameStyle = $true

# set paragraph mark formatting: range covering just the pilcrow (End to End, or Paragraph.Range.End-? )
$markRange = $d.Range($p2.Range.End - 1, $p2.Range.End)
Write-Output "markRange.Text = [$($markRange.Text)]"
$markRange.HighlightColorIndex = 3
$markRange.Font.TextColor.ObjectThemeColor = 13

$runRange = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:r>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:highlight w:val="cyan"/>
    </w:rPr>
    <w:t xml:space="preserve">This was written by Luiza Sartori on the </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:highlight w:val="cyan"/>
    </w:rPr>
    <w:t>22/03/2023</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
</w:p>
</w:body>
</w:document>
'@
$runRange.InsertXML($xml)
Write-Output "done"
